# =====================================================================
# Edit script for results_MPC_example_2018_01_01__02_00_00.xlsx
#
# Summary of the change (per commit message / diff):
#  1. A new worksheet "storage_content__Building1" is inserted right
#     before "costs__Building1". It reports, per timestep, the storage
#     content of the shStorage and electricalStorage components (this
#     data used to live as an extra "storage_content" column on the
#     shSourceBus__Building1 sheet).
#  2. The "storage_content" column (E) is removed from
#     shSourceBus__Building1 (dimension A1:E26 -> A1:D26).
#  3. On env_impacts__Building1 the two storage rows (6 & 7) swap
#     order: shStorage now comes before electricalStorage.
#  4. On capStorages__Building1 the two storage rows (2 & 3) swap
#     order: shStorage now comes before electricalStorage.
#  5. On spaceHeatingBus__Building1, row 22 (timestamp 43101.91666...)
#     gets updated flow values for columns B and C (column D is kept).
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "storage_content__Building1" sheet before "costs"
# ---------------------------------------------------------------------
$costsSheet = $wb.Worksheets.Item("costs__Building1")
$storageSheet = $wb.Worksheets.Add($costsSheet)
$storageSheet.Name = "storage_content__Building1"

# Header row (bold, centered, bordered - same look as the other sheets)
$storageSheet.Cells.Item(1, 2).Value = "shStorage__B001_storage_content"
$storageSheet.Cells.Item(1, 3).Value = "electricalStorage__B001_storage_content"
$headerRange = $storageSheet.Range("B1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows: Date/time (col A), shStorage content (col B),
# electricalStorage content (col C)
$storageData = @(
    @(43101.08333333334, 26.64435025980561, 81.8440278108984),
    @(43101.125, 13.36468957599664, 81.20497465508444),
    @(43101.16666666666, 0, 80.47625372485189),
    @(43101.20833333334, 0, 76.09120489044363),
    @(43101.25, 0, 68.87869781831),
    @(43101.29166666666, 0, 62.44957082839697),
    @(43101.33333333334, 0, 55.58070722114262),
    @(43101.375, 0, 51.92463975228887),
    @(43101.41666666666, 0, 54.97232508470272),
    @(43101.45833333334, 0, 54.97232508470272),
    @(43101.5, 52.24448642575736, 54.97232508470272),
    @(43101.54166666666, 52.325, 60.15545109839242),
    @(43101.58333333334, 52.325, 57.69631328583431),
    @(43101.625, 52.325, 53.96085001796482),
    @(43101.66666666666, 44.3433705550596, 52.33473706098807),
    @(43101.70833333334, 35.31286976631884, 50.19596629703458),
    @(43101.75, 25.83785739473363, 47.8635908819183),
    @(43101.79166666666, 18.49078851382524, 37.80593307494156),
    @(43101.83333333334, 11.81135031878672, 23.72043473773225),
    @(43101.875, 52.325, 6.377749323933808),
    @(43101.91666666666, 48.95808314575903, 2.812009967441861),
    @(43101.95833333334, 37.91711655679155, 1.865880398837209),
    @(43102, 25.9466953121862, 1.194684386046512),
    @(43102.04166666666, 13.09308468147299, 0.5555980069767442),
    @(43102.08333333334, 0, 0)
)

$r = 2
foreach ($row in $storageData) {
    $storageSheet.Cells.Item($r, 1).Value = $row[0]
    $storageSheet.Cells.Item($r, 2).Value = $row[1]
    $storageSheet.Cells.Item($r, 3).Value = $row[2]
    $r++
}

$lastRow = $r - 1
$dateRange = $storageSheet.Range("A2:A" + $lastRow)
$dateRange.Font.Bold = $true
$dateRange.HorizontalAlignment = -4108
$dateRange.VerticalAlignment = -4160
$dateRange.Borders.LineStyle = 1
$dateRange.NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# 2. Remove the "storage_content" column (E) from shSourceBus
# ---------------------------------------------------------------------
$shSourceBus = $wb.Worksheets.Item("shSourceBus__Building1")
$shSourceBus.Range("E1:E26").Delete()

# ---------------------------------------------------------------------
# 3. Swap the storage rows on env_impacts__Building1 (rows 6 & 7)
# ---------------------------------------------------------------------
$envImpacts = $wb.Worksheets.Item("env_impacts__Building1")
$envImpacts.Cells.Item(6, 1).Value = "shStorage__Building1"
$envImpacts.Cells.Item(6, 2).Value = 22.37625
$envImpacts.Cells.Item(7, 1).Value = "electricalStorage__Building1"
$envImpacts.Cells.Item(7, 2).Value = 2757.066666666667

# ---------------------------------------------------------------------
# 4. Swap the storage rows on capStorages__Building1 (rows 2 & 3)
# ---------------------------------------------------------------------
$capStorages = $wb.Worksheets.Item("capStorages__Building1")
$capStorages.Cells.Item(2, 1).Value = "shStorage__Building1"
$capStorages.Cells.Item(2, 2).Value = 4500
$capStorages.Cells.Item(3, 1).Value = "electricalStorage__Building1"
$capStorages.Cells.Item(3, 2).Value = 200

# ---------------------------------------------------------------------
# 5. Update row 22 flows on spaceHeatingBus__Building1
# ---------------------------------------------------------------------
$spaceHeatingBus = $wb.Worksheets.Item("spaceHeatingBus__Building1")
$spaceHeatingBus.Cells.Item(22, 2).Value = 6.093526585699433
$spaceHeatingBus.Cells.Item(22, 3).Value = 3.232733646300566

Write-Output "Edit complete"
